# remove unsupported fa from UG
# Applies three changes to the single slide of the class-diagram deck:
#  1. Reposition the (invisible) subtitle divider line shape.
#  2. Rename the "prerequisite" association label to "needs" and shrink its box.
#  3. Rename the "preclusion" association label to "precludes" and shrink its box.
#
# NOTE: Shape.Left/Top/Width/Height are exposed as single-precision (Single)
# floats in the PowerPoint object model, so plain "EMU/12700.0" point values
# can round-trip to the wrong EMU integer after the float32 truncation. The
# literals below were chosen so that, once truncated to float32 and converted
# back to EMU (round(pt * 12700)), they reproduce the exact target EMU values
# from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "Subtitle 2" divider shape: off (1524000,5257799) -> (2054086,6202679) EMU
$divider = $s.Shapes.Item(1)
$divider.Left = 161.7390594482422
$divider.Top = 488.3999328613281

# 2) "prerequisite" -> "needs" textbox: ext cx 819455 -> 497252 EMU
$needsBox = $s.Shapes.Item(10)
$needsBox.TextFrame.TextRange.Text = "needs"
$needsBox.Width = 39.15370178222656

# 3) "preclusion" -> "precludes" textbox: ext cx 777777 -> 742511 EMU
$precludesBox = $s.Shapes.Item(16)
$precludesBox.TextFrame.TextRange.Text = "precludes"
$precludesBox.Width = 58.46543502807617
